# Insert two new data rows (for a new reporting date: 44468) right above the
# current row 273, pushing the existing data (rows 273:369) down to rows
# 275:371. The new rows replicate the structure of the rows that used to be
# at 273/274 (same Mercado/Region/Categoria/Calidad/etc.) but carry their own
# Fecha / Volumen / Precio values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 273:274 -- everything at/after row 273 shifts down by 2.
$ws.Rows("273:274").Insert()

# --- New row 273 ("Primera") ---
$ws.Range("A273").Value = 3
$ws.Range("B273").Value = "Femacal de La Calera"
$ws.Range("C273").Value = "Coquimbo"
$ws.Range("D273").Value = 44468
$ws.Range("E273").Value = 5
$ws.Range("F273").Value = 100112008
$ws.Range("G273").Value = "Coliflor"
$ws.Range("H273").Value = "Sin especificar"
$ws.Range("I273").Value = "Primera"
$ws.Range("J273").Value = 2200
$ws.Range("K273").Value = 650
$ws.Range("L273").Value = 700
$ws.Range("M273").Value = 677
$ws.Range("N273").Value = "`$/unidad"
$ws.Range("O273").Value = "Provincia de Quillota"
$ws.Range("P273").Value = 677
$ws.Range("Q273").Value = 1
$ws.Range("R273").Value = "Hortaliza"

# --- New row 274 ("Segunda") ---
$ws.Range("A274").Value = 3
$ws.Range("B274").Value = "Femacal de La Calera"
$ws.Range("C274").Value = "Coquimbo"
$ws.Range("D274").Value = 44468
$ws.Range("E274").Value = 5
$ws.Range("F274").Value = 100112008
$ws.Range("G274").Value = "Coliflor"
$ws.Range("H274").Value = "Sin especificar"
$ws.Range("I274").Value = "Segunda"
$ws.Range("J274").Value = 1100
$ws.Range("K274").Value = 550
$ws.Range("L274").Value = 550
$ws.Range("M274").Value = 550
$ws.Range("N274").Value = "`$/unidad"
$ws.Range("O274").Value = "Provincia de Quillota"
$ws.Range("P274").Value = 550
$ws.Range("Q274").Value = 1
$ws.Range("R274").Value = "Hortaliza"
